# Added medical equipment indicator mappings for TB for Budget Mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndicatorMappingDB")

# --- Fix typo "Actiivity" -> "Activity" for all existing Indicator Type
# --- Name values in column H (rows 2-342 are the "Activity" indicator type).
for ($r = 2; $r -le 342; $r++) {
    $ws.Cells.Item($r, 8).Value = "Activity"
}

# --- Append new rows: Medical equipment indicator mappings for TB
# --- (Module = Infrastructure, Indicator Type = Medical equipment).
$newRows = @(
    @{ Row = 415; IndicatorId = 1; IndicatorName = "Health post" },
    @{ Row = 416; IndicatorId = 2; IndicatorName = "Health centre" },
    @{ Row = 417; IndicatorId = 3; IndicatorName = "District/General hospital" },
    @{ Row = 418; IndicatorId = 4; IndicatorName = "National/Regional/Provincial hospital" },
    @{ Row = 419; IndicatorId = 6; IndicatorName = "Prehospital emergency" },
    @{ Row = 420; IndicatorId = 7; IndicatorName = "Free-standing general outpatient clinic" },
    @{ Row = 421; IndicatorId = 8; IndicatorName = "Free-standing specialized outpatient clinic" }
)

# First pass: every column except the Indicator Name (J), in row order. This
# establishes "Infrastructure" / "Medical equipment" as new shared strings.
foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = 2                                         # A: Costing Modes
    $ws.Cells.Item($r, 2).Value = "TB Costing"                              # B: Costing Mode Names
    $ws.Cells.Item($r, 3).Value = 1                                         # C: Budget ID
    $ws.Cells.Item($r, 4).Value = "Global Fund"                            # D: Budget Name
    $ws.Cells.Item($r, 5).Value = 20                                        # E: Module ID
    $ws.Cells.Item($r, 6).Value = "Infrastructure"                         # F: Module Name
    $ws.Cells.Item($r, 7).Value = 4                                         # G: Indicator Type ID
    $ws.Cells.Item($r, 8).Value = "Medical equipment"                     # H: Indicator Type Name
    $ws.Cells.Item($r, 9).Value = $nr.IndicatorId                          # I: Indicator ID
    $ws.Cells.Item($r, 17).Value = 56                                       # Q: Budget Category ID
    $ws.Cells.Item($r, 18).Value = "DS,TB screening and diagnosis"        # R: Budget Category Name
}

# Second pass: Indicator Name (J) values, written in the same order the
# original author first typed/introduced them (Health post, Prehospital
# emergency, Health centre, District/General hospital, National/.../
# Provincial hospital, Free-standing general/specialized outpatient clinic)
# so new shared strings land at the same table positions as the source file.
$indicatorNameOrder = @(415, 419, 416, 417, 418, 420, 421)
foreach ($r in $indicatorNameOrder) {
    $nr = $newRows | Where-Object { $_.Row -eq $r }
    $ws.Cells.Item($r, 10).Value = $nr.IndicatorName                      # J: Indicator Name
}
